# The two occurrence records on rows 19 and 20 were swapped: the data that
# used to be reported on row 19 (Djupsvart brunbagge / Melandrya dubia,
# observed 2023-08-17) now belongs on row 20, and the data that used to be
# on row 20 (Vedtrappmossa / Crossocalyx hellerianus, observed 2023-08-18)
# now belongs on row 19. The "Aktivitet" (M) and "Publik kommentar" (AC)
# cells travel with the row-19 record, so they move from row 19 to row 20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the current ("before") values of both rows -------------------
$A19 = $ws.Range("A19").Value2
$B19 = $ws.Range("B19").Value2
$D19 = $ws.Range("D19").Value2
$E19 = $ws.Range("E19").Value2
$F19 = $ws.Range("F19").Value2
$G19 = $ws.Range("G19").Value2
$H19 = $ws.Range("H19").Value2
$Q19 = $ws.Range("Q19").Value2
$R19 = $ws.Range("R19").Value2
$Y19 = $ws.Range("Y19").Value2
$AA19 = $ws.Range("AA19").Value2
$AC19 = $ws.Range("AC19").Value2

$A20 = $ws.Range("A20").Value2
$B20 = $ws.Range("B20").Value2
$D20 = $ws.Range("D20").Value2
$E20 = $ws.Range("E20").Value2
$F20 = $ws.Range("F20").Value2
$G20 = $ws.Range("G20").Value2
$H20 = $ws.Range("H20").Value2
$Q20 = $ws.Range("Q20").Value2
$R20 = $ws.Range("R20").Value2
$Y20 = $ws.Range("Y20").Value2
$AA20 = $ws.Range("AA20").Value2

# --- write row 20's old data into row 19 -----------------------------------
$ws.Range("A19").Value2 = $A20
$ws.Range("B19").Value2 = $B20
$ws.Range("D19").Value2 = $D20
$ws.Range("E19").Value2 = $E20
$ws.Range("F19").Value2 = $F20
$ws.Range("G19").Value2 = $G20
$ws.Range("H19").Value2 = $H20
$ws.Range("Q19").Value2 = $Q20
$ws.Range("R19").Value2 = $R20

# Y19 / AA19 hold plain text dates ("2023-08-18"), not real date values;
# force text so Excel doesn't auto-convert the string to a date serial,
# then restore the default "Normal" style so no stray formatting is left
# behind on the cell.
$ws.Range("Y19").NumberFormat = "@"
$ws.Range("Y19").Value2 = $Y20
$ws.Range("Y19").Style = "Normal"

$ws.Range("AA19").NumberFormat = "@"
$ws.Range("AA19").Value2 = $AA20
$ws.Range("AA19").Style = "Normal"

# Row 19 no longer carries the "Aktivitet" / "Publik kommentar" notes -
# those move to row 20 below.
$ws.Range("M19").ClearContents()
$ws.Range("AC19").ClearContents()

# --- write row 19's old data into row 20 -----------------------------------
$ws.Range("A20").Value2 = $A19
$ws.Range("B20").Value2 = $B19
$ws.Range("D20").Value2 = $D19
$ws.Range("E20").Value2 = $E19
$ws.Range("F20").Value2 = $F19
$ws.Range("G20").Value2 = $G19
$ws.Range("H20").Value2 = $H19
$ws.Range("Q20").Value2 = $Q19
$ws.Range("R20").Value2 = $R19

$ws.Range("Y20").NumberFormat = "@"
$ws.Range("Y20").Value2 = $Y19
$ws.Range("Y20").Style = "Normal"

$ws.Range("AA20").NumberFormat = "@"
$ws.Range("AA20").Value2 = $AA19
$ws.Range("AA20").Style = "Normal"

# Row 20 now gets the "Aktivitet" (empty) and "Publik kommentar" cells that
# used to sit on row 19.
$ws.Range("M20").NumberFormat = "@"
$ws.Range("M20").Value2 = ""
$ws.Range("M20").Style = "Normal"

$ws.Range("AC20").NumberFormat = "@"
$ws.Range("AC20").Value2 = $AC19
$ws.Range("AC20").Style = "Normal"
